# Fruta / hortaliza, semanal
# Insert two new weekly report rows right before the existing row 117,
# pushing the previously-existing rows 117-148 down to 119-150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 117 (shifts 117..148 -> 119..150)
$ws.Rows.Item(117).Insert()
$ws.Rows.Item(117).Insert()

# New row 117: Brócoli "Primera" entry for the latest week
$ws.Cells.Item(117, 1).Value = 11
$ws.Cells.Item(117, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(117, 3).Value = "Bíobío"
$ws.Cells.Item(117, 4).Value = 44463
$ws.Cells.Item(117, 5).Value = 8
$ws.Cells.Item(117, 6).Value = 100112023
$ws.Cells.Item(117, 7).Value = "Brócoli"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 2000
$ws.Cells.Item(117, 11).Value = 700
$ws.Cells.Item(117, 12).Value = 800
$ws.Cells.Item(117, 13).Value = 750
$ws.Cells.Item(117, 14).Value = "$/unidad"
$ws.Cells.Item(117, 15).Value = "Región Metropolitana"
$ws.Cells.Item(117, 16).Value = 750
$ws.Cells.Item(117, 17).Value = 1
$ws.Cells.Item(117, 18).Value = "Hortaliza"

# New row 118: Brócoli "Segunda" entry for the latest week
$ws.Cells.Item(118, 1).Value = 11
$ws.Cells.Item(118, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(118, 3).Value = "Bíobío"
$ws.Cells.Item(118, 4).Value = 44463
$ws.Cells.Item(118, 5).Value = 8
$ws.Cells.Item(118, 6).Value = 100112023
$ws.Cells.Item(118, 7).Value = "Brócoli"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Segunda"
$ws.Cells.Item(118, 10).Value = 1000
$ws.Cells.Item(118, 11).Value = 600
$ws.Cells.Item(118, 12).Value = 600
$ws.Cells.Item(118, 13).Value = 600
$ws.Cells.Item(118, 14).Value = "$/unidad"
$ws.Cells.Item(118, 15).Value = "Región Metropolitana"
$ws.Cells.Item(118, 16).Value = 600
$ws.Cells.Item(118, 17).Value = 1
$ws.Cells.Item(118, 18).Value = "Hortaliza"
